{"js": "// 1) Merge the two \"Accuracy Rate\" runs that were split by the stray\n//    \"_GoBack\" bookmark back into a single contiguous run/sentence.\n// 2) Remove the \"F-score\" paragraph together with the extra blank\n//    paragraph that trailed it, leaving one blank paragraph behind.\n// 3) Move the \"_GoBack\" bookmark out of the \"Accuracy Rate\" paragraph and\n//    into that surviving blank paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// Locate the paragraphs we need by their (unique) leading text.\nlet accuracyIndex = -1;\nlet fScoreIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (accuracyIndex === -1 && t.indexOf(\"Accuracy Rate\") === 0) {\n    accuracyIndex = i;\n  }\n  if (fScoreIndex === -1 && t.indexOf(\"F-score\") === 0) {\n    fScoreIndex = i;\n  }\n}\nif (accuracyIndex === -1) {\n  throw new Error(\"Could not find the Accuracy Rate paragraph\");\n}\nif (fScoreIndex === -1) {\n  throw new Error(\"Could not find the F-score paragraph\");\n}\n\nconst accuracyParagraph = paragraphs.items[accuracyIndex];\nconst survivingBlank = paragraphs.items[fScoreIndex - 1]; // blank paragraph kept\nconst fScoreParagraph = paragraphs.items[fScoreIndex];\nconst blankAfterFScore = paragraphs.items[fScoreIndex + 1]; // blank paragraph dropped\n\n// Remove the bookmark from its current mid-sentence position; we'll\n// reinsert it afterwards.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Collapse the \"Accuracy Rate\" paragraph's explanatory text (previously\n// split across two runs) back into one run by replacing everything after\n// the bold \"Accuracy Rate\" label with its own (now contiguous) text.\nconst boldLabel = accuracyParagraph.search(\"Accuracy Rate\", { matchCase: true });\nboldLabel.load(\"items\");\nawait context.sync();\nconst afterLabel = boldLabel.items[0]\n  .getRange(\"After\")\n  .expandTo(accuracyParagraph.getRange(\"End\"));\nafterLabel.load(\"text\");\nawait context.sync();\nafterLabel.insertText(afterLabel.text, \"Replace\");\nawait context.sync();\n\n// Delete the F-score paragraph and the blank paragraph that used to\n// follow it, keeping only the blank paragraph that preceded it.\nfScoreParagraph.delete();\nblankAfterFScore.delete();\nawait context.sync();\n\n// Put the bookmark back, now living alone in the surviving blank\n// paragraph right after \"Accuracy Rate\".\nsurvivingBlank.getRange().insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# 1) Merge the two \"Accuracy Rate\" runs that were split by the stray\n#    \"_GoBack\" bookmark back into a single contiguous run/sentence.\n# 2) Remove the \"F-score\" paragraph together with the extra blank\n#    paragraph that trailed it, leaving one blank paragraph behind.\n# 3) Move the \"_GoBack\" bookmark out of the \"Accuracy Rate\" paragraph and\n#    into that surviving blank paragraph.\n\n$d = $word.ActiveDocument\n\n# --- Locate the paragraphs we need to touch, by their distinctive text. ---\n$accuracyIndex = -1\n$fScoreIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($accuracyIndex -eq -1 -and $t.StartsWith(\"Accuracy Rate\")) {\n        $accuracyIndex = $i\n    }\n    if ($fScoreIndex -eq -1 -and $t.StartsWith(\"F-score\")) {\n        $fScoreIndex = $i\n    }\n}\nif ($accuracyIndex -eq -1) { throw \"Could not find the Accuracy Rate paragraph\" }\nif ($fScoreIndex -eq -1) { throw \"Could not find the F-score paragraph\" }\n\n# --- Step 1: collapse the two runs of the \"Accuracy Rate\" paragraph into one. ---\n# The sentence was split in two by the \"_GoBack\" bookmark; a Find/Replace over\n# the full (bookmark-spanning) sentence rewrites it as a single contiguous run\n# and - as a side effect - removes the bookmark that used to sit in the middle.\n$accuracyParagraph = $d.Paragraphs.Item($accuracyIndex)\n$explanation = \": The accuracy rate is a measure of how well the classifier predicts the correct class labels. It is calculated by dividing the total number of correct predictions (sum of the diagonal elements in the confusion matrix) by the total number of samples. In this case, the accuracy rate is approximately 0.805, indicating that the classifier correctly predicts the class labels for about 80.5% of the test set samples.\"\n\n$find = $accuracyParagraph.Range.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.MatchCase = $true\n$find.Text = $explanation\n$find.Replacement.Text = $explanation\n$find.Execute([ref]$find.Text, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$find.Replacement.Text, [ref]2) | Out-Null\n\n# --- Step 2: delete the \"F-score\" paragraph and the blank paragraph that ---\n# --- trailed it, keeping the blank paragraph that precedes it.          ---\n$fScoreParagraph = $d.Paragraphs.Item($fScoreIndex)\n$blankAfterStart = $fScoreParagraph.Range.Start\n$blankAfterEnd = $d.Paragraphs.Item($fScoreIndex + 1).Range.End\n$toDelete = $d.Range($blankAfterStart, $blankAfterEnd)\n$toDelete.Delete()\n\n# --- Step 3: put the bookmark back in the surviving blank paragraph. ---\n$survivingBlank = $d.Paragraphs.Item($fScoreIndex - 1)\n$bookmarkRange = $survivingBlank.Range\n$bookmarkRange.Collapse(1)  # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange) | Out-Null\n"}
